# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Numeric-looking price strings are written as text ("@" number format, then
# Style reset to "Normal" to avoid leaving a stray format behind) so they
# keep matching the sheet's original inline-string representation instead
# of being auto-coerced to floats by Excel's Range.Value parser.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.673.51"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "3.849.40"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "457.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +13.71%  "
$ws.Range("E7").Value = "  +3.18%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.744"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.155"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000319"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").Value = "4.446.03"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.61%  "
$ws.Range("D16").Value = "3.857.88"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.56%  "
$ws.Range("E19").Value = "  +7.92%  "
$ws.Range("D20").Value = "67.692.95"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "429.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.78"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.51%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.35%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "739.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.134"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.36%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0477"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.354"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +14.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +18.86%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0688"
$ws.Range("E43").Value = "  -7.40%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.27%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.140"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.24%  "
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.30%  "
$ws.Range("E48").Value = "  +5.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.10%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.59%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
